# ============================================================
# Excel COM-interop script implementing the commit:
# "refactor: Integrate cascade events into unified JSON structure"
#
# 1) Trims Tradeoff_Relationships / Uncertainty_Weights / Test_Constants
#    down to their new, simplified column sets + refreshed values.
# 2) Appends 6 brand-new constants worksheets at the end of the workbook:
#    Metric_Ranges, Game_Flow_Constants, Probability_Constants,
#    Threshold_Constants, Storyteller_Constants, Technical_Constants.
# ============================================================

$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd([string]$name) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newWs = $wb.Worksheets.Add($null, $last)
    $newWs.Name = $name
    return $newWs
}

function Write-Table($ws, $rows) {
    for ($r = 0; $r -lt $rows.Count; $r++) {
        $row = $rows[$r]
        for ($c = 0; $c -lt $row.Count; $c++) {
            $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
        }
    }
}

# Reusable "header style" source (bold, centered, bordered - style index 1
# in the original workbook) so new header rows match the existing look
# instead of Excel's default cell style.
$headerStyleSrc = $wb.Worksheets.Item("Game_Constants").Range("A1")

# ---- Update existing sheet: Tradeoff_Relationships ----
$rows3new = @(
    @("Source_Metric","Target_Metric","Impact_Factor","Description"),
    @("MONEY","HAPPINESS",-0.5,"MONEY 상승 시 HAPPINESS 하락"),
    @("MONEY","STAFF_FATIGUE",-0.5,"MONEY 상승 시 STAFF_FATIGUE 하락"),
    @("REPUTATION","MONEY",-0.5,"REPUTATION 상승 시 MONEY 하락"),
    @("REPUTATION","STAFF_FATIGUE",-0.5,"REPUTATION 상승 시 STAFF_FATIGUE 하락"),
    @("HAPPINESS","SUFFERING",-0.5,"HAPPINESS 상승 시 SUFFERING 하락"),
    @("SUFFERING","HAPPINESS",-0.5,"SUFFERING 상승 시 HAPPINESS 하락"),
    @("INVENTORY","MONEY",-0.5,"INVENTORY 상승 시 MONEY 하락"),
    @("STAFF_FATIGUE","REPUTATION",-0.5,"STAFF_FATIGUE 상승 시 REPUTATION 하락"),
    @("STAFF_FATIGUE","FACILITY",-0.5,"STAFF_FATIGUE 상승 시 FACILITY 하락"),
    @("FACILITY","MONEY",-0.5,"FACILITY 상승 시 MONEY 하락"),
    @("DEMAND","INVENTORY",-0.5,"DEMAND 상승 시 INVENTORY 하락"),
    @("DEMAND","STAFF_FATIGUE",-0.5,"DEMAND 상승 시 STAFF_FATIGUE 하락")
)
$ws3 = $wb.Worksheets.Item("Tradeoff_Relationships")
$oldLastCol3 = $ws3.UsedRange.Columns.Count
$oldLastRow3 = $ws3.UsedRange.Rows.Count
for ($col = $oldLastCol3; $col -gt 4; $col--) { $ws3.Columns.Item($col).Delete() }
for ($row = $oldLastRow3; $row -gt 13; $row--) { $ws3.Rows.Item($row).Delete() }
Write-Table $ws3 $rows3new
$ws3.Range("A2:D13").ClearFormats()
$headerStyleSrc.Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# ---- Update existing sheet: Uncertainty_Weights ----
$rows4new = @(
    @("Metric_Name","Weight","Description"),
    @("MONEY",0.3,"MONEY의 불확실성 가중치"),
    @("REPUTATION",0.25,"REPUTATION의 불확실성 가중치"),
    @("HAPPINESS",-0.1,"HAPPINESS의 불확실성 가중치"),
    @("SUFFERING",0.2,"SUFFERING의 불확실성 가중치"),
    @("INVENTORY",0.05,"INVENTORY의 불확실성 가중치"),
    @("STAFF_FATIGUE",0.15,"STAFF_FATIGUE의 불확실성 가중치"),
    @("FACILITY",-0.2,"FACILITY의 불확실성 가중치"),
    @("DEMAND",0.1,"DEMAND의 불확실성 가중치")
)
$ws4 = $wb.Worksheets.Item("Uncertainty_Weights")
$oldLastCol4 = $ws4.UsedRange.Columns.Count
$oldLastRow4 = $ws4.UsedRange.Rows.Count
for ($col = $oldLastCol4; $col -gt 3; $col--) { $ws4.Columns.Item($col).Delete() }
for ($row = $oldLastRow4; $row -gt 9; $row--) { $ws4.Rows.Item($row).Delete() }
Write-Table $ws4 $rows4new
$ws4.Range("A2:C9").ClearFormats()
$headerStyleSrc.Copy()
$ws4.Range("A1:C1").PasteSpecial(-4122)

# ---- Update existing sheet: Test_Constants ----
$rows9new = @(
    @("Key","Value","Type","Category","Description"),
    @("TEST_MIN_CASCADE_EVENTS",3,"int","testing","최소 연쇄 효과 메시지 수"),
    @("TEST_EXPECTED_EVENTS",2,"int","testing","예상 이벤트 수"),
    @("TEST_METRICS_HISTORY_LENGTH",5,"int","testing","메트릭 히스토리 길이")
)
$ws9 = $wb.Worksheets.Item("Test_Constants")
$oldLastCol9 = $ws9.UsedRange.Columns.Count
$oldLastRow9 = $ws9.UsedRange.Rows.Count
for ($col = $oldLastCol9; $col -gt 5; $col--) { $ws9.Columns.Item($col).Delete() }
for ($row = $oldLastRow9; $row -gt 4; $row--) { $ws9.Rows.Item($row).Delete() }
Write-Table $ws9 $rows9new
$ws9.Range("A2:E4").ClearFormats()
$headerStyleSrc.Copy()
$ws9.Range("A1:E1").PasteSpecial(-4122)

# ---- New sheet: Metric_Ranges ----
$rows12 = @(
    @("Metric_Name","Min_Value","Max_Value","Default_Value","Description"),
    @("MONEY",0,"inf",10000,"MONEY의 허용 범위 및 기본값"),
    @("REPUTATION",0,100,50,"REPUTATION의 허용 범위 및 기본값"),
    @("HAPPINESS",0,100,50,"HAPPINESS의 허용 범위 및 기본값"),
    @("SUFFERING",0,100,20,"SUFFERING의 허용 범위 및 기본값"),
    @("INVENTORY",0,"inf",100,"INVENTORY의 허용 범위 및 기본값"),
    @("STAFF_FATIGUE",0,100,30,"STAFF_FATIGUE의 허용 범위 및 기본값"),
    @("FACILITY",0,100,80,"FACILITY의 허용 범위 및 기본값"),
    @("DEMAND",0,"inf",60,"DEMAND의 허용 범위 및 기본값")
)
$ws12 = Add-SheetAtEnd "Metric_Ranges"
Write-Table $ws12 $rows12
$headerStyleSrc.Copy()
$ws12.Range("A1:E1").PasteSpecial(-4122)

# ---- New sheet: Game_Flow_Constants ----
$rows13 = @(
    @("Key","Value","Type","Category","Description"),
    @("MAX_ACTIONS_PER_DAY",3,"int","game_flow","하루 최대 행동 횟수"),
    @("DEFAULT_GAME_LENGTH",30,"int","game_flow","기본 게임 길이(일)"),
    @("DEFAULT_TOTAL_DAYS",730,"int","game_flow","기본 게임 총 일수"),
    @("DEFAULT_COOLDOWN_DAYS",5,"int","game_flow","기본 쿨다운 일수")
)
$ws13 = Add-SheetAtEnd "Game_Flow_Constants"
Write-Table $ws13 $rows13
$headerStyleSrc.Copy()
$ws13.Range("A1:E1").PasteSpecial(-4122)

# ---- New sheet: Probability_Constants ----
$rows14 = @(
    @("Key","Value","Type","Category","Description"),
    @("PROBABILITY_LOW_THRESHOLD",0.3,"float","probability","낮은 확률 임계값"),
    @("PROBABILITY_HIGH_THRESHOLD",0.7,"float","probability","높은 확률 임계값"),
    @("DEFAULT_PROBABILITY",0.8,"float","probability","기본 확률값"),
    @("DEFAULT_SEVERITY",0.5,"float","probability","기본 심각도")
)
$ws14 = Add-SheetAtEnd "Probability_Constants"
Write-Table $ws14 $rows14
$headerStyleSrc.Copy()
$ws14.Range("A1:E1").PasteSpecial(-4122)

# ---- New sheet: Threshold_Constants ----
$rows15 = @(
    @("Key","Value","Type","Category","Description"),
    @("MONEY_LOW_THRESHOLD",3000,"int","thresholds","자금 부족 기준"),
    @("MONEY_HIGH_THRESHOLD",15000,"int","thresholds","자금 풍부 기준"),
    @("REPUTATION_LOW_THRESHOLD",30,"int","thresholds","평판 위험 기준"),
    @("REPUTATION_HIGH_THRESHOLD",70,"int","thresholds","평판 우수 기준"),
    @("HAPPINESS_LOW_THRESHOLD",30,"int","thresholds","행복 위험 기준"),
    @("HAPPINESS_HIGH_THRESHOLD",70,"int","thresholds","행복 우수 기준"),
    @("REPUTATION_BASELINE",50,"int","thresholds","평판 기준점")
)
$ws15 = Add-SheetAtEnd "Threshold_Constants"
Write-Table $ws15 $rows15
$headerStyleSrc.Copy()
$ws15.Range("A1:E1").PasteSpecial(-4122)

# ---- New sheet: Storyteller_Constants ----
$rows16 = @(
    @("Key","Value","Type","Category","Description"),
    @("MIN_METRICS_HISTORY_FOR_TREND",2,"int","storyteller","추세 분석을 위한 최소 히스토리 개수"),
    @("RECENT_HISTORY_WINDOW",3,"int","storyteller","최근 히스토리 분석 윈도우 크기"),
    @("MINIMUM_TREND_POINTS",2,"int","storyteller","트렌드 분석에 필요한 최소 데이터 포인트"),
    @("SITUATION_POSITIVE_THRESHOLD",0.6,"float","storyteller","긍정적 상황 판단 임계값"),
    @("SITUATION_NEGATIVE_THRESHOLD",0.4,"float","storyteller","부정적 상황 판단 임계값"),
    @("TRADEOFF_BALANCE_THRESHOLD",0.5,"float","storyteller","트레이드오프 불균형 감지 임계값"),
    @("GAME_PROGRESSION_MID_POINT",0.5,"float","storyteller","게임 진행도 중간점"),
    @("PATTERN_SCORE_TOLERANCE",0.1,"float","storyteller","패턴 점수 허용 오차"),
    @("COMPLEXITY_BONUS_MULTIPLIER",0.1,"float","storyteller","복잡성 보너스 배수")
)
$ws16 = Add-SheetAtEnd "Storyteller_Constants"
Write-Table $ws16 $rows16
$headerStyleSrc.Copy()
$ws16.Range("A1:E1").PasteSpecial(-4122)

# ---- New sheet: Technical_Constants ----
$rows17 = @(
    @("Key","Value","Type","Category","Description"),
    @("FLOAT_EPSILON",0.001,"float","technical","부동소수점 비교 오차 허용 범위"),
    @("SCORE_THRESHOLD_HIGH",0.7,"float","technical","높은 점수 임계값"),
    @("SCORE_THRESHOLD_MEDIUM",0.5,"float","technical","중간 점수 임계값")
)
$ws17 = Add-SheetAtEnd "Technical_Constants"
Write-Table $ws17 $rows17
$headerStyleSrc.Copy()
$ws17.Range("A1:E1").PasteSpecial(-4122)

Write-Host "Done. Worksheets:" $wb.Worksheets.Count
